$wb = $excel.ActiveWorkbook

# Sheet ALC, row 26 (item id 1963)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H26").Value = 3750
$ws.Range("I26").Value = 2500
$ws.Range("J26").Value = 5000
$ws.Range("K26").Value = 2500
$ws.Range("L26").Value = 5000
$ws.Range("M26").Value = -2156
$ws.Range("N26").Value = -5688

# Sheet ALC, row 96 (item id 19894)
$ws.Range("H96").Value = 446.84616
$ws.Range("I96").Value = 673.4286
$ws.Range("J96").Value = 182.5
$ws.Range("K96").Value = 2020.2858
$ws.Range("L96").Value = 547.5
$ws.Range("M96").Value = -647.2857999999999
$ws.Range("N96").Value = -3293.5

# Sheet ALC, row 106 (item id 19903)
$ws.Range("H106").Value = 6674.5
$ws.Range("I106").Value = 6699.25
$ws.Range("J106").Value = 6625
$ws.Range("K106").Value = 6699.25
$ws.Range("L106").Value = 6625
$ws.Range("M106").Value = -6068.25
$ws.Range("N106").Value = -7887

# Sheet ALC, row 132 (item id 44049)
$ws.Range("H132").Value = 3463.182
$ws.Range("I132").Value = 3173.2104
$ws.Range("K132").Value = 9519.6312
$ws.Range("M132").Value = -6989.6312

# Sheet ALC, row 135 (item id 44047)
$ws.Range("H135").Value = 1959.3704
$ws.Range("J135").Value = 4100.8335
$ws.Range("L135").Value = 36907.5015
$ws.Range("N135").Value = -41977.5015

# Sheet ALC, row 141 (item id 44161)
$ws.Range("H141").Value = 1947.9
$ws.Range("I141").Value = 1964.3334
$ws.Range("K141").Value = 5893.0002
$ws.Range("M141").Value = -713.0002000000004

# Sheet ARM, row 2 (item id 27713)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1720.7368
$ws.Range("I2").Value = 1734.258
$ws.Range("J2").Value = 1660.8572
$ws.Range("K2").Value = 1734.258
$ws.Range("L2").Value = 1660.8572
$ws.Range("M2").Value = -1621.258
$ws.Range("N2").Value = -1886.8572

# Sheet ARM, row 32 (item id 44147)
$ws.Range("H32").Value = 10886.079
$ws.Range("I32").Value = 6893.3228
$ws.Range("K32").Value = 6893.3228
$ws.Range("M32").Value = -6606.3228

# Sheet ARM, row 33 (item id 3352)
$ws.Range("H33").Value = 54984.8
$ws.Range("I33").Value = 52308.332
$ws.Range("J33").Value = 58999.5
$ws.Range("K33").Value = 52308.332
$ws.Range("L33").Value = 58999.5
$ws.Range("M33").Value = -51979.332
$ws.Range("N33").Value = -59657.5

# Sheet ARM, row 45 (item id 27714)
$ws.Range("H45").Value = 1705.4375
$ws.Range("I45").Value = 1103.0769
$ws.Range("J45").Value = 4315.6665
$ws.Range("K45").Value = 1103.0769
$ws.Range("L45").Value = 4315.6665
$ws.Range("M45").Value = -726.0769
$ws.Range("N45").Value = -5069.6665

# Sheet ARM, row 116 (item id 27713)
$ws.Range("H116").Value = 1720.7368
$ws.Range("I116").Value = 1734.258
$ws.Range("J116").Value = 1660.8572
$ws.Range("K116").Value = 1734.258
$ws.Range("L116").Value = 1660.8572
$ws.Range("M116").Value = 559.742
$ws.Range("N116").Value = -6248.8572

# Sheet ARM, row 132 (item id 43997)
$ws.Range("H132").Value = 12699.333
$ws.Range("I132").Value = 13434.444
$ws.Range("J132").Value = 6083.3335
$ws.Range("K132").Value = 40303.33199999999
$ws.Range("L132").Value = 18250.0005
$ws.Range("M132").Value = -37773.33199999999
$ws.Range("N132").Value = -23310.0005

# Sheet BSM, row 3 (item id 27713)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1720.7368
$ws.Range("I3").Value = 1734.258
$ws.Range("J3").Value = 1660.8572
$ws.Range("K3").Value = 1734.258
$ws.Range("L3").Value = 1660.8572
$ws.Range("M3").Value = -1620.258
$ws.Range("N3").Value = -1888.8572

# Sheet BSM, row 22 (item id 5092)
$ws.Range("H22").Value = 291
$ws.Range("I22").Value = 294
$ws.Range("J22").Value = 285
$ws.Range("K22").Value = 294
$ws.Range("L22").Value = 285
$ws.Range("M22").Value = -121
$ws.Range("N22").Value = -631

# Sheet BSM, row 86 (item id 12526)
$ws.Range("H86").Value = 52727730
$ws.Range("I86").Value = 62557730
$ws.Range("J86").Value = 301060
$ws.Range("K86").Value = 62557730
$ws.Range("L86").Value = 301060
$ws.Range("M86").Value = -62556607
$ws.Range("N86").Value = -303306

# Sheet BSM, row 89 (item id 12526)
$ws.Range("H89").Value = 52727730
$ws.Range("I89").Value = 62557730
$ws.Range("J89").Value = 301060
$ws.Range("K89").Value = 312788650
$ws.Range("L89").Value = 1505300
$ws.Range("M89").Value = -312783034
$ws.Range("N89").Value = -1516532

# Sheet BSM, row 107 (item id 27706)
$ws.Range("H107").Value = 2282.3333
$ws.Range("I107").Value = 2282.3333
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 2282.3333
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -362.3332999999998
$ws.Range("N107").ClearContents()

# Sheet CRP, row 16 (item id 27691)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 1139
$ws.Range("I16").Value = 703.5
$ws.Range("J16").Value = 2010
$ws.Range("K16").Value = 703.5
$ws.Range("L16").Value = 2010
$ws.Range("M16").Value = -416.5
$ws.Range("N16").Value = -2584

# Sheet CRP, row 31 (item id 44023)
$ws.Range("H31").Value = 2259.0344
$ws.Range("I31").Value = 1697.9333
$ws.Range("J31").Value = 2860.2144
$ws.Range("K31").Value = 1697.9333
$ws.Range("L31").Value = 2860.2144
$ws.Range("M31").Value = -1402.9333
$ws.Range("N31").Value = -3450.2144

# Sheet CRP, row 34 (item id 44023)
$ws.Range("H34").Value = 2259.0344
$ws.Range("I34").Value = 1697.9333
$ws.Range("J34").Value = 2860.2144
$ws.Range("K34").Value = 1697.9333
$ws.Range("L34").Value = 2860.2144
$ws.Range("M34").Value = -1495.9333
$ws.Range("N34").Value = -3264.2144

# Sheet CRP, row 58 (item id 44021)
$ws.Range("H58").Value = 3951.6538
$ws.Range("I58").Value = 3211.6956
$ws.Range("K58").Value = 3211.6956
$ws.Range("M58").Value = -3008.6956

# Sheet CRP, row 113 (item id 27691)
$ws.Range("H113").Value = 1139
$ws.Range("I113").Value = 703.5
$ws.Range("J113").Value = 2010
$ws.Range("K113").Value = 703.5
$ws.Range("L113").Value = 2010
$ws.Range("M113").Value = 1466.5
$ws.Range("N113").Value = -6350

# Sheet CRP, row 133 (item id 43328)
$ws.Range("H133").Value = 150000
$ws.Range("J133").Value = 150000
$ws.Range("L133").Value = 150000
$ws.Range("N133").Value = -155060

# Sheet CRP, row 136 (item id 44021)
$ws.Range("H136").Value = 3951.6538
$ws.Range("I136").Value = 3211.6956
$ws.Range("K136").Value = 9635.086800000001
$ws.Range("M136").Value = -7085.086800000001

# Sheet CUL, row 63 (item id 12866)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H63").Value = 14374.75
$ws.Range("I63").Value = 3750
$ws.Range("K63").Value = 11250
$ws.Range("M63").Value = -10501

# Sheet CUL, row 66 (item id 12866)
$ws.Range("H66").Value = 14374.75
$ws.Range("I66").Value = 3750
$ws.Range("K66").Value = 33750
$ws.Range("M66").Value = -30006

# Sheet CUL, row 114 (item id 27865)
$ws.Range("H114").Value = 2930.077
$ws.Range("J114").Value = 2541.8572
$ws.Range("L114").Value = 7625.571599999999
$ws.Range("N114").Value = -14133.5716

# Sheet CUL, row 117 (item id 27870)
$ws.Range("H117").Value = 687.2857
$ws.Range("I117").Value = 458
$ws.Range("J117").Value = 993
$ws.Range("K117").Value = 1374
$ws.Range("L117").Value = 2979
$ws.Range("M117").Value = 2068
$ws.Range("N117").Value = -9863

# Sheet CUL, row 136 (item id 44093)
$ws.Range("H136").Value = 5628.375
$ws.Range("I136").Value = 5837.8335
$ws.Range("J136").Value = 5000
$ws.Range("K136").Value = 17513.5005
$ws.Range("L136").Value = 15000
$ws.Range("M136").Value = -12413.5005
$ws.Range("N136").Value = -25200

# Sheet GSM, row 15 (item id 12018)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H15").Value = 40535.332
$ws.Range("I15").Value = 42738
$ws.Range("J15").Value = 38332.668
$ws.Range("K15").Value = 42738
$ws.Range("L15").Value = 38332.668
$ws.Range("M15").Value = -42450
$ws.Range("N15").Value = -38908.668

# Sheet GSM, row 70 (item id 14146)
$ws.Range("H70").Value = 6815
$ws.Range("I70").Value = 6172
$ws.Range("J70").Value = 7200.8
$ws.Range("K70").Value = 6172
$ws.Range("L70").Value = 7200.8
$ws.Range("M70").Value = -5902
$ws.Range("N70").Value = -7740.8

# Sheet GSM, row 73 (item id 14146)
$ws.Range("H73").Value = 6815
$ws.Range("I73").Value = 6172
$ws.Range("J73").Value = 7200.8
$ws.Range("K73").Value = 6172
$ws.Range("L73").Value = 7200.8
$ws.Range("M73").Value = -5236
$ws.Range("N73").Value = -9072.799999999999

# Sheet GSM, row 81 (item id 12018)
$ws.Range("H81").Value = 40535.332
$ws.Range("I81").Value = 42738
$ws.Range("J81").Value = 38332.668
$ws.Range("K81").Value = 42738
$ws.Range("L81").Value = 38332.668
$ws.Range("M81").Value = -41740
$ws.Range("N81").Value = -40328.668

# Sheet GSM, row 84 (item id 12018)
$ws.Range("H84").Value = 40535.332
$ws.Range("I84").Value = 42738
$ws.Range("J84").Value = 38332.668
$ws.Range("K84").Value = 128214
$ws.Range("L84").Value = 114998.004
$ws.Range("M84").Value = -123222
$ws.Range("N84").Value = -124982.004

# Sheet GSM, row 102 (item id 36169)
$ws.Range("H102").Value = 28430.25
$ws.Range("I102").Value = 34321.156
$ws.Range("J102").Value = 12721.167
$ws.Range("K102").Value = 34321.156
$ws.Range("L102").Value = 12721.167
$ws.Range("M102").Value = -32699.156
$ws.Range("N102").Value = -15965.167

# Sheet GSM, row 113 (item id 27710)
$ws.Range("H113").Value = 35723700
$ws.Range("I113").Value = 47629640
$ws.Range("J113").Value = 5883.857
$ws.Range("K113").Value = 47629640
$ws.Range("L113").Value = 5883.857
$ws.Range("M113").Value = -47627470
$ws.Range("N113").Value = -10223.857

# Sheet LTW, row 19 (item id 2229)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H19").Value = 6281.8
$ws.Range("I19").Value = 1900
$ws.Range("J19").Value = 9203
$ws.Range("K19").Value = 1900
$ws.Range("L19").Value = 9203
$ws.Range("M19").Value = -1730
$ws.Range("N19").Value = -9543

# Sheet LTW, row 40 (item id 36248)
$ws.Range("H40").Value = 5197.6875
$ws.Range("I40").Value = 5391.2856
$ws.Range("K40").Value = 5391.2856
$ws.Range("M40").Value = -5255.2856

# Sheet LTW, row 132 (item id 44058)
$ws.Range("H132").Value = 2819.818
$ws.Range("I132").Value = 2791
$ws.Range("J132").Value = 2949.5
$ws.Range("K132").Value = 8373
$ws.Range("L132").Value = 8848.5
$ws.Range("M132").Value = -5843
$ws.Range("N132").Value = -13908.5

# Sheet WVR, row 107 (item id 27746)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 1150.4546
$ws.Range("I107").Value = 1150.4546
$ws.Range("J107").Value = 0
$ws.Range("K107").Value = 3451.3638
$ws.Range("L107").Value = 0
$ws.Range("M107").Value = -1531.3638
$ws.Range("N107").ClearContents()

# Sheet WVR, row 108 (item id 25661)
$ws.Range("H108").Value = 0
$ws.Range("J108").Value = 0
$ws.Range("L108").Value = 0
$ws.Range("N108").ClearContents()

# Sheet WVR, row 132 (item id 44029)
$ws.Range("H132").Value = 9714
$ws.Range("I132").Value = 8241.615
$ws.Range("K132").Value = 24724.845
$ws.Range("M132").Value = -22194.845

